$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date format is applied to column D for rows 275-283 (matches style s="2" / numFmt 165)
$ws.Range('D275:D283').NumberFormat = 'YYYY-MM-DD HH:MM:SS'

# Row 275
$ws.Cells.Item(275, 1).Value = 7
$ws.Cells.Item(275, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(275, 3).Value = 'Ñuble'
$ws.Cells.Item(275, 4).Value = 44939
$ws.Cells.Item(275, 5).Value = 16
$ws.Cells.Item(275, 6).Value = 100112027
$ws.Cells.Item(275, 7).Value = 'Melón'
$ws.Cells.Item(275, 8).Value = 'Calameño'
$ws.Cells.Item(275, 9).Value = 'Extra'
$ws.Cells.Item(275, 10).Value = 500
$ws.Cells.Item(275, 11).Value = 1200
$ws.Cells.Item(275, 12).Value = 1200
$ws.Cells.Item(275, 13).Value = 1200
$ws.Cells.Item(275, 14).Value = '$/unidad'
$ws.Cells.Item(275, 15).Value = 'Región del Maule'
$ws.Cells.Item(275, 16).Value = 1200
$ws.Cells.Item(275, 17).Value = 1
$ws.Cells.Item(275, 18).Value = 'Hortaliza'

# Row 276
$ws.Cells.Item(276, 1).Value = 7
$ws.Cells.Item(276, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(276, 3).Value = 'Ñuble'
$ws.Cells.Item(276, 4).Value = 44939
$ws.Cells.Item(276, 5).Value = 16
$ws.Cells.Item(276, 6).Value = 100112027
$ws.Cells.Item(276, 7).Value = 'Melón'
$ws.Cells.Item(276, 8).Value = 'Calameño'
$ws.Cells.Item(276, 9).Value = 'Primera'
$ws.Cells.Item(276, 10).Value = 500
$ws.Cells.Item(276, 11).Value = 1000
$ws.Cells.Item(276, 12).Value = 1000
$ws.Cells.Item(276, 13).Value = 1000
$ws.Cells.Item(276, 14).Value = '$/unidad'
$ws.Cells.Item(276, 15).Value = 'Región del Maule'
$ws.Cells.Item(276, 16).Value = 1000
$ws.Cells.Item(276, 17).Value = 1
$ws.Cells.Item(276, 18).Value = 'Hortaliza'

# Row 277
$ws.Cells.Item(277, 1).Value = 7
$ws.Cells.Item(277, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(277, 3).Value = 'Ñuble'
$ws.Cells.Item(277, 4).Value = 44939
$ws.Cells.Item(277, 5).Value = 16
$ws.Cells.Item(277, 6).Value = 100112027
$ws.Cells.Item(277, 7).Value = 'Melón'
$ws.Cells.Item(277, 8).Value = 'Calameño'
$ws.Cells.Item(277, 9).Value = 'Segunda'
$ws.Cells.Item(277, 10).Value = 500
$ws.Cells.Item(277, 11).Value = 800
$ws.Cells.Item(277, 12).Value = 800
$ws.Cells.Item(277, 13).Value = 800
$ws.Cells.Item(277, 14).Value = '$/unidad'
$ws.Cells.Item(277, 15).Value = 'Región del Maule'
$ws.Cells.Item(277, 16).Value = 800
$ws.Cells.Item(277, 17).Value = 1
$ws.Cells.Item(277, 18).Value = 'Hortaliza'

# Row 278
$ws.Cells.Item(278, 1).Value = 7
$ws.Cells.Item(278, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(278, 3).Value = 'Ñuble'
$ws.Cells.Item(278, 4).Value = 44939
$ws.Cells.Item(278, 5).Value = 16
$ws.Cells.Item(278, 6).Value = 100112027
$ws.Cells.Item(278, 7).Value = 'Melón'
$ws.Cells.Item(278, 8).Value = 'Tuna'
$ws.Cells.Item(278, 9).Value = 'Extra'
$ws.Cells.Item(278, 10).Value = 500
$ws.Cells.Item(278, 11).Value = 1200
$ws.Cells.Item(278, 12).Value = 1200
$ws.Cells.Item(278, 13).Value = 1200
$ws.Cells.Item(278, 14).Value = '$/unidad'
$ws.Cells.Item(278, 15).Value = 'Región del Maule'
$ws.Cells.Item(278, 16).Value = 1200
$ws.Cells.Item(278, 17).Value = 1
$ws.Cells.Item(278, 18).Value = 'Hortaliza'

# Row 279
$ws.Cells.Item(279, 1).Value = 7
$ws.Cells.Item(279, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(279, 3).Value = 'Ñuble'
$ws.Cells.Item(279, 4).Value = 44939
$ws.Cells.Item(279, 5).Value = 16
$ws.Cells.Item(279, 6).Value = 100112027
$ws.Cells.Item(279, 7).Value = 'Melón'
$ws.Cells.Item(279, 8).Value = 'Tuna'
$ws.Cells.Item(279, 9).Value = 'Primera'
$ws.Cells.Item(279, 10).Value = 500
$ws.Cells.Item(279, 11).Value = 1000
$ws.Cells.Item(279, 12).Value = 1000
$ws.Cells.Item(279, 13).Value = 1000
$ws.Cells.Item(279, 14).Value = '$/unidad'
$ws.Cells.Item(279, 15).Value = 'Región del Maule'
$ws.Cells.Item(279, 16).Value = 1000
$ws.Cells.Item(279, 17).Value = 1
$ws.Cells.Item(279, 18).Value = 'Hortaliza'

# Row 280
$ws.Cells.Item(280, 1).Value = 7
$ws.Cells.Item(280, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(280, 3).Value = 'Ñuble'
$ws.Cells.Item(280, 4).Value = 44939
$ws.Cells.Item(280, 5).Value = 16
$ws.Cells.Item(280, 6).Value = 100112027
$ws.Cells.Item(280, 7).Value = 'Melón'
$ws.Cells.Item(280, 8).Value = 'Tuna'
$ws.Cells.Item(280, 9).Value = 'Segunda'
$ws.Cells.Item(280, 10).Value = 500
$ws.Cells.Item(280, 11).Value = 800
$ws.Cells.Item(280, 12).Value = 800
$ws.Cells.Item(280, 13).Value = 800
$ws.Cells.Item(280, 14).Value = '$/unidad'
$ws.Cells.Item(280, 15).Value = 'Región del Maule'
$ws.Cells.Item(280, 16).Value = 800
$ws.Cells.Item(280, 17).Value = 1
$ws.Cells.Item(280, 18).Value = 'Hortaliza'

# Row 281
$ws.Cells.Item(281, 1).Value = 7
$ws.Cells.Item(281, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(281, 3).Value = 'Ñuble'
$ws.Cells.Item(281, 4).Value = 44209
$ws.Cells.Item(281, 5).Value = 16
$ws.Cells.Item(281, 6).Value = 100112027
$ws.Cells.Item(281, 7).Value = 'Melón'
$ws.Cells.Item(281, 8).Value = 'Calameño'
$ws.Cells.Item(281, 9).Value = 'Extra'
$ws.Cells.Item(281, 10).Value = 1200
$ws.Cells.Item(281, 11).Value = 950
$ws.Cells.Item(281, 12).Value = 1000
$ws.Cells.Item(281, 13).Value = 975
$ws.Cells.Item(281, 14).Value = '$/unidad'
$ws.Cells.Item(281, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(281, 16).Value = 975
$ws.Cells.Item(281, 17).Value = 1
$ws.Cells.Item(281, 18).Value = 'Hortaliza'

# Row 282
$ws.Cells.Item(282, 1).Value = 7
$ws.Cells.Item(282, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(282, 3).Value = 'Ñuble'
$ws.Cells.Item(282, 4).Value = 44209
$ws.Cells.Item(282, 5).Value = 16
$ws.Cells.Item(282, 6).Value = 100112027
$ws.Cells.Item(282, 7).Value = 'Melón'
$ws.Cells.Item(282, 8).Value = 'Calameño'
$ws.Cells.Item(282, 9).Value = 'Primera'
$ws.Cells.Item(282, 10).Value = 2100
$ws.Cells.Item(282, 11).Value = 750
$ws.Cells.Item(282, 12).Value = 800
$ws.Cells.Item(282, 13).Value = 779
$ws.Cells.Item(282, 14).Value = '$/unidad'
$ws.Cells.Item(282, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(282, 16).Value = 779
$ws.Cells.Item(282, 17).Value = 1
$ws.Cells.Item(282, 18).Value = 'Hortaliza'

# Row 283
$ws.Cells.Item(283, 1).Value = 7
$ws.Cells.Item(283, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(283, 3).Value = 'Ñuble'
$ws.Cells.Item(283, 4).Value = 44209
$ws.Cells.Item(283, 5).Value = 16
$ws.Cells.Item(283, 6).Value = 100112027
$ws.Cells.Item(283, 7).Value = 'Melón'
$ws.Cells.Item(283, 8).Value = 'Calameño'
$ws.Cells.Item(283, 9).Value = 'Segunda'
$ws.Cells.Item(283, 10).Value = 1800
$ws.Cells.Item(283, 11).Value = 550
$ws.Cells.Item(283, 12).Value = 600
$ws.Cells.Item(283, 13).Value = 578
$ws.Cells.Item(283, 14).Value = '$/unidad'
$ws.Cells.Item(283, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(283, 16).Value = 578
$ws.Cells.Item(283, 17).Value = 1
$ws.Cells.Item(283, 18).Value = 'Hortaliza'

